{"js": "// Find the sentence body (without the leading \"- \" and the trailing\n// period) and recolor it green, which splits the original single run\n// into \"- \" / body / \".\" pieces in the underlying OOXML.\nconst body = context.document.body;\nconst target =\n  \"Has task owner\\u2019s basic details, task description and \" +\n  \"\\u201CFurther details\\u201D if owner decided to add it\";\n\nconst results = body.search(target, { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nconst sentenceBody = results.items[0];\nsentenceBody.font.color = \"#00B050\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Build the target sentence (without the leading \"- \" and the trailing\n# period) using explicit Unicode code points for the curly\n# apostrophe/quotes so Find matches the text exactly.\n$apos = [char]0x2019\n$lq   = [char]0x201C\n$rq   = [char]0x201D\n$target = \"Has task owner\" + $apos + \"s basic details, task description and \" + $lq + \"Further details\" + $rq + \" if owner decided to add it\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $target\n$find.Execute() | Out-Null\n\nif ($find.Found) {\n    # Recoloring just this inner span splits the original run into\n    # \"- \" / body / \".\" pieces, leaving the dash and period uncolored.\n    $find.Parent.Font.Color = 5287936\n}\n"}
